$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$colOrder = @('B','C','D','E','F','I','J','K','L','M','N')
$data = @(
    ,@(2, 1.02, 1.039482991695807, 1.049141092787546, 1.038006646078971, 1.056966896595004, 1.038537309730384, 1.044574764889604, 1.051898978263679, 1.040795908726703, 1.059703157406805, 1.046058180870148)
    ,@(3, 1.02, 1.040533216866516, 1.050136796849661, 1.038901613526406, 1.05804283016121, 1.038743283158534, 1.045269710415986, 1.052706394639569, 1.041500575928293, 1.060592167654349, 1.046754113298889)
    ,@(4, 1.02, 1.041212688302554, 1.050781251546125, 1.039481011611792, 1.058739208209519, 1.038874570948507, 1.04571871624163, 1.053228401842152, 1.041956204945406, 1.061166995973935, 1.047203756764306)
    ,@(5, 1.02, 1.041498315567358, 1.051052220045382, 1.039724660486286, 1.059032007638121, 1.038929287360746, 1.045907317492677, 1.053447746767585, 1.042147670147702, 1.061408552743239, 1.047392625850763)
    ,@(6, 1.02, 1.041546272323548, 1.051097719168365, 1.03976557428882, 1.059081172410417, 1.038938446506078, 1.045938975034433, 1.053484569467373, 1.04217981322304, 1.061449105252133, 1.047424328349859)
    ,@(7, 1.02, 1.041216504956453, 1.050784872084789, 1.039484266985929, 1.058743120444788, 1.038875303946553, 1.045721236972997, 1.053231333160257, 1.041958763632491, 1.061170224064816, 1.047206281075401)
    ,@(8, 1.02, 1.039837939827371, 1.049477560713133, 1.038309043259821, 1.057330476537733, 1.03860733139762, 1.0448097636144, 1.052171940393039, 1.041034124025581, 1.060003689508887, 1.046293513320088)
    ,@(9, 1.02, 1.03740800350786, 1.047175210398829, 1.036240425979571, 1.054842586045736, 1.038119898509646, 1.043198509344289, 1.050301759029009, 1.039402219656237, 1.057944897689818, 1.044679970884099)
    ,@(10, 1.02, 1.035787545856876, 1.045641202028071, 1.034862907347321, 1.053184925634463, 1.037784722918254, 1.042120910794823, 1.049052703726812, 1.038312573545274, 1.056570227562026, 1.043600842021091)
    ,@(11, 1.02, 1.035085747773233, 1.044977173845906, 1.034266801695413, 1.052467363381182, 1.037637167811695, 1.041653487639179, 1.048511313477833, 1.037840342420042, 1.055974474315494, 1.043132755070951)
    ,@(12, 1.02, 1.034825048882975, 1.044730555514266, 1.034045437285276, 1.052200861215898, 1.037581995645171, 1.041479743370125, 1.048310135827434, 1.037664873752097, 1.055753108169362, 1.042958764065107)
    ,@(13, 1.02, 1.034880970584225, 1.044783454565319, 1.034092918173692, 1.05225802536857, 1.037593846713608, 1.041517017647671, 1.048353292789874, 1.037702515134797, 1.055800595436545, 1.042996091276402)
    ,@(14, 1.02, 1.035064198728694, 1.044956787658493, 1.034248502504616, 1.052445333567992, 1.037632614679846, 1.041639128382167, 1.04849468572663, 1.037825839368666, 1.05595617768849, 1.043118375422147)
    ,@(15, 1.02, 1.035177088994934, 1.045063588049002, 1.034344370545795, 1.052560744625701, 1.037656452725553, 1.041714348606067, 1.048581791823217, 1.037901815421032, 1.056052026865159, 1.043193702467394)
    ,@(16, 1.02, 1.035834119125974, 1.045685275774856, 1.034902476719182, 1.053232552422078, 1.037794464658382, 1.042151914947577, 1.049088622622189, 1.038343905405898, 1.056609754960669, 1.043631890203297)
    ,@(17, 1.02, 1.036246221876182, 1.046075299716156, 1.035252661006315, 1.053654017444881, 1.037880387658469, 1.042426170378415, 1.049406399394526, 1.038621107885147, 1.056959466024459, 1.043906535108306)
    ,@(18, 1.02, 1.036486581997962, 1.04630281410201, 1.035456953158204, 1.053899871668448, 1.037930271401052, 1.042586060206363, 1.049591700999833, 1.038782756104402, 1.057163397166206, 1.044066651998149)
    ,@(19, 1.02, 1.03656853649065, 1.046380394020105, 1.035526617527732, 1.053983705134053, 1.037947240823331, 1.042640565150195, 1.049654875200382, 1.03883786729908, 1.057232924033623, 1.044121234345252)
    ,@(20, 1.02, 1.036202008440794, 1.046033451766704, 1.035215085863511, 1.053608796084331, 1.037871193103631, 1.042396753508406, 1.049372310353833, 1.038591370745876, 1.056921950446581, 1.043877076462967)
    ,@(21, 1.02, 1.035010243157902, 1.044905744517795, 1.03420268525482, 1.05239017507723, 1.037621208516223, 1.0416031731953, 1.04845305127195, 1.037789525138175, 1.055910364714467, 1.043082369174791)
    ,@(22, 1.02, 1.034260818057464, 1.044196891363825, 1.033566471641514, 1.051624167168028, 1.037461929541651, 1.041103509060375, 1.047874606802403, 1.037285020121827, 1.055273895925932, 1.042581995459482)
    ,@(23, 1.02, 1.034658113521344, 1.044572650624188, 1.033903709805638, 1.05203022468192, 1.037546565675692, 1.041368457623573, 1.04818129566606, 1.037552501167807, 1.055611342206391, 1.042847320280031)
    ,@(24, 1.02, 1.036221986627302, 1.046052360991747, 1.035232064343122, 1.05362922960618, 1.037875348449203, 1.042410045968859, 1.049387713885168, 1.03860480780042, 1.056938902273893, 1.043890387800238)
    ,@(25, 1.02, 1.03803628766208, 1.047770267258149, 1.036774940316406, 1.055485601813737, 1.038247714449618, 1.043615662626397, 1.05078564655299, 1.039824408788658, 1.058477523397512, 1.045097716571717)
)

foreach ($rowData in $data) {
    $r = $rowData[0]
    for ($i = 0; $i -lt $colOrder.Length; $i++) {
        $col = $colOrder[$i]
        $val = $rowData[$i + 1]
        $ws.Range("$col$r").Value = $val
    }
}
